# CreateAcount.xlsx - add PhoneNumber / CurrentEmail columns, clean up
# leftover blank-placeholder values, and populate sample input data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (K, L) for "PhoneNumber" / "CurrentEmail" before
# the existing "NameMesEr" column; everything from the old K onward shifts
# right by two columns.
$ws.Range("K1:L1").EntireColumn.Insert()

# New header cells.
$ws.Range("K1").Value = "PhoneNumber"
$ws.Range("L1").Value = "CurrentEmail"

# Row 2 ("All field is blank") - drop the stray whitespace-only values that
# used to sit in FirstName/LastName.
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

# Row 3 ("Only firstname and lastname are blank") - same cleanup, plus the
# new sample phone/email values.
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("L3").Value = "dinhhuy131@gmail.com"
$ws.Range("K3").Value = "'0979155626"
$ws.Range("K3").NumberFormat = "@"

# Row 4 ("Only firstname is blank") - clear stray FirstName value, keep
# LastName, add new sample data.
$ws.Range("B4").ClearContents()
$ws.Range("L4").Value = "dinhhuy131@gmail.com"
$ws.Range("K4").Value = "'0979155627"
$ws.Range("K4").NumberFormat = "@"

# Row 5 ("Only lastname is blank") - FirstName stays, just add new sample
# data.
$ws.Range("L5").Value = "dinhhuy131@gmail.com"
$ws.Range("K5").Value = "'0979155628"
$ws.Range("K5").NumberFormat = "@"

# Match the author's final selection.
$ws.Range("B14").Select() | Out-Null
